# Append a new login record (test6) to the "login" sheet,
# mirroring the existing username/password rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("login")

$ws.Range("A7").Value = "test6@gmail.com"
$ws.Range("B7").Value = "test6"
